$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'ECs'
$ws.Range("B2").Value = 'Hp'
$ws.Range("C2").Value = 'Itgam'
$ws.Range("D2").Value = 'FAPs'
$ws.Range("E2").Value = [double]"1"
$ws.Range("F2").Value = [double]"0.3333333333333333"
$ws.Range("G2").Value = [double]"0.8909086666666667"
$ws.Range("H2").Value = [double]"2.672726"
$ws.Range("I2").Value = [double]"0.03117145726690431"
$ws.Range("J2").Value = [double]"0.0311714572669043"
$ws.Range("K2").Value = [double]"1"
$ws.Range("L2").Value = [double]"0.3333333333333333"
$ws.Range("M2").Value = [double]"0.142723"
$ws.Range("N2").Value = [double]"0.428169"
$ws.Range("O2").Value = [double]"0.0009642800942465787"
$ws.Range("P2").Value = [double]"0.0009642800942465787"
$ws.Range("Q2").Value = [double]"0.1271531576326667"
$ws.Range("R2").Value = [double]"1.144378418694"
$ws.Range("S2").Value = [double]"3.005801575113369E-05"
$ws.Range("T2").Value = [double]"3.005801575113368E-05"

# Row 3
$ws.Range("A3").Value = 'ECs'
$ws.Range("B3").Value = 'Hp'
$ws.Range("C3").Value = 'Itgam'
$ws.Range("D3").Value = 'Inflammatory-Mac'
$ws.Range("E3").Value = [double]"1"
$ws.Range("F3").Value = [double]"0.3333333333333333"
$ws.Range("G3").Value = [double]"0.8909086666666667"
$ws.Range("H3").Value = [double]"2.672726"
$ws.Range("I3").Value = [double]"0.03117145726690431"
$ws.Range("J3").Value = [double]"0.0311714572669043"
$ws.Range("K3").Value = [double]"3"
$ws.Range("L3").Value = [double]"1"
$ws.Range("M3").Value = [double]"86.42780700000002"
$ws.Range("N3").Value = [double]"259.283421"
$ws.Range("O3").Value = [double]"0.5839326098770704"
$ws.Range("P3").Value = [double]"0.5839326098770704"
$ws.Range("Q3").Value = [double]"76.99928229729402"
$ws.Range("R3").Value = [double]"692.993540675646"
$ws.Range("S3").Value = [double]"0.01820203039553501"
$ws.Range("T3").Value = [double]"0.018202030395535"

# Row 4
$ws.Range("A4").Value = 'ECs'
$ws.Range("B4").Value = 'Hp'
$ws.Range("C4").Value = 'Itgam'
$ws.Range("D4").Value = 'MuSCs'
$ws.Range("E4").Value = [double]"1"
$ws.Range("F4").Value = [double]"0.3333333333333333"
$ws.Range("G4").Value = [double]"0.8909086666666667"
$ws.Range("H4").Value = [double]"2.672726"
$ws.Range("I4").Value = [double]"0.03117145726690431"
$ws.Range("J4").Value = [double]"0.0311714572669043"
$ws.Range("K4").Value = [double]"1"
$ws.Range("L4").Value = [double]"0.3333333333333333"
$ws.Range("M4").Value = [double]"0.006361333333333333"
$ws.Range("N4").Value = [double]"0.019084"
$ws.Range("O4").Value = [double]"4.297910712499435E-05"
$ws.Range("P4").Value = [double]"4.297910712499435E-05"
$ws.Range("Q4").Value = [double]"0.005667366998222223"
$ws.Range("R4").Value = [double]"0.051006302984"
$ws.Range("S4").Value = [double]"1.339721401116464E-06"
$ws.Range("T4").Value = [double]"1.339721401116464E-06"

# Row 5
$ws.Range("A5").Value = 'ECs'
$ws.Range("B5").Value = 'Hp'
$ws.Range("C5").Value = 'Itgam'
$ws.Range("D5").Value = 'Resolving-Mac'
$ws.Range("E5").Value = [double]"1"
$ws.Range("F5").Value = [double]"0.3333333333333333"
$ws.Range("G5").Value = [double]"0.8909086666666667"
$ws.Range("H5").Value = [double]"2.672726"
$ws.Range("I5").Value = [double]"0.03117145726690431"
$ws.Range("J5").Value = [double]"0.0311714572669043"
$ws.Range("K5").Value = [double]"3"
$ws.Range("L5").Value = [double]"1"
$ws.Range("M5").Value = [double]"61.43300833333333"
$ws.Range("N5").Value = [double]"184.299025"
$ws.Range("O5").Value = [double]"0.415060130921558"
$ws.Range("P5").Value = [double]"0.415060130921558"
$ws.Range("Q5").Value = [double]"54.73119954357222"
$ws.Range("R5").Value = [double]"492.58079589215"
$ws.Range("S5").Value = [double]"0.01293802913421705"
$ws.Range("T5").Value = [double]"0.01293802913421705"

# Row 6
$ws.Range("A6").Value = 'FAPs'
$ws.Range("B6").Value = 'Hp'
$ws.Range("C6").Value = 'Itgam'
$ws.Range("D6").Value = 'FAPs'
$ws.Range("E6").Value = [double]"3"
$ws.Range("F6").Value = [double]"1"
$ws.Range("G6").Value = [double]"15.02703733333333"
$ws.Range("H6").Value = [double]"45.081112"
$ws.Range("I6").Value = [double]"0.5257717986252713"
$ws.Range("J6").Value = [double]"0.5257717986252713"
$ws.Range("K6").Value = [double]"1"
$ws.Range("L6").Value = [double]"0.3333333333333333"
$ws.Range("M6").Value = [double]"0.142723"
$ws.Range("N6").Value = [double]"0.428169"
$ws.Range("O6").Value = [double]"0.0009642800942465787"
$ws.Range("P6").Value = [double]"0.0009642800942465787"
$ws.Range("Q6").Value = [double]"2.144703849325334"
$ws.Range("R6").Value = [double]"19.302334643928"
$ws.Range("S6").Value = [double]"0.0005069912795305698"
$ws.Range("T6").Value = [double]"0.0005069912795305698"

# Row 7
$ws.Range("A7").Value = 'FAPs'
$ws.Range("B7").Value = 'Hp'
$ws.Range("C7").Value = 'Itgam'
$ws.Range("D7").Value = 'Inflammatory-Mac'
$ws.Range("E7").Value = [double]"3"
$ws.Range("F7").Value = [double]"1"
$ws.Range("G7").Value = [double]"15.02703733333333"
$ws.Range("H7").Value = [double]"45.081112"
$ws.Range("I7").Value = [double]"0.5257717986252713"
$ws.Range("J7").Value = [double]"0.5257717986252713"
$ws.Range("K7").Value = [double]"3"
$ws.Range("L7").Value = [double]"1"
$ws.Range("M7").Value = [double]"86.42780700000002"
$ws.Range("N7").Value = [double]"259.283421"
$ws.Range("O7").Value = [double]"0.5839326098770704"
$ws.Range("P7").Value = [double]"0.5839326098770704"
$ws.Range("Q7").Value = [double]"1298.753882427128"
$ws.Range("R7").Value = [double]"11688.78494184415"
$ws.Range("S7").Value = [double]"0.3070152985710162"
$ws.Range("T7").Value = [double]"0.3070152985710162"

# Row 8
$ws.Range("A8").Value = 'FAPs'
$ws.Range("B8").Value = 'Hp'
$ws.Range("C8").Value = 'Itgam'
$ws.Range("D8").Value = 'MuSCs'
$ws.Range("E8").Value = [double]"3"
$ws.Range("F8").Value = [double]"1"
$ws.Range("G8").Value = [double]"15.02703733333333"
$ws.Range("H8").Value = [double]"45.081112"
$ws.Range("I8").Value = [double]"0.5257717986252713"
$ws.Range("J8").Value = [double]"0.5257717986252713"
$ws.Range("K8").Value = [double]"1"
$ws.Range("L8").Value = [double]"0.3333333333333333"
$ws.Range("M8").Value = [double]"0.006361333333333333"
$ws.Range("N8").Value = [double]"0.019084"
$ws.Range("O8").Value = [double]"4.297910712499435E-05"
$ws.Range("P8").Value = [double]"4.297910712499435E-05"
$ws.Range("Q8").Value = [double]"0.09559199348977779"
$ws.Range("R8").Value = [double]"0.8603279414080001"
$ws.Range("S8").Value = [double]"2.259720245641649E-05"
$ws.Range("T8").Value = [double]"2.259720245641649E-05"

# Row 9
$ws.Range("A9").Value = 'FAPs'
$ws.Range("B9").Value = 'Hp'
$ws.Range("C9").Value = 'Itgam'
$ws.Range("D9").Value = 'Resolving-Mac'
$ws.Range("E9").Value = [double]"3"
$ws.Range("F9").Value = [double]"1"
$ws.Range("G9").Value = [double]"15.02703733333333"
$ws.Range("H9").Value = [double]"45.081112"
$ws.Range("I9").Value = [double]"0.5257717986252713"
$ws.Range("J9").Value = [double]"0.5257717986252713"
$ws.Range("K9").Value = [double]"3"
$ws.Range("L9").Value = [double]"1"
$ws.Range("M9").Value = [double]"61.43300833333333"
$ws.Range("N9").Value = [double]"184.299025"
$ws.Range("O9").Value = [double]"0.415060130921558"
$ws.Range("P9").Value = [double]"0.415060130921558"
$ws.Range("Q9").Value = [double]"923.1561097239778"
$ws.Range("R9").Value = [double]"8308.404987515802"
$ws.Range("S9").Value = [double]"0.2182269115722681"
$ws.Range("T9").Value = [double]"0.2182269115722681"

# Row 10
$ws.Range("A10").Value = 'Inflammatory-Mac'
$ws.Range("B10").Value = 'Hp'
$ws.Range("C10").Value = 'Itgam'
$ws.Range("D10").Value = 'FAPs'
$ws.Range("E10").Value = [double]"3"
$ws.Range("F10").Value = [double]"1"
$ws.Range("G10").Value = [double]"12.47652966666667"
$ws.Range("H10").Value = [double]"37.429589"
$ws.Range("I10").Value = [double]"0.436533649177391"
$ws.Range("J10").Value = [double]"0.436533649177391"
$ws.Range("K10").Value = [double]"1"
$ws.Range("L10").Value = [double]"0.3333333333333333"
$ws.Range("M10").Value = [double]"0.142723"
$ws.Range("N10").Value = [double]"0.428169"
$ws.Range("O10").Value = [double]"0.0009642800942465787"
$ws.Range("P10").Value = [double]"0.0009642800942465787"
$ws.Range("Q10").Value = [double]"1.780687743615667"
$ws.Range("R10").Value = [double]"16.026189692541"
$ws.Range("S10").Value = [double]"0.0004209407083705775"
$ws.Range("T10").Value = [double]"0.0004209407083705775"

# Row 11
$ws.Range("A11").Value = 'Inflammatory-Mac'
$ws.Range("B11").Value = 'Hp'
$ws.Range("C11").Value = 'Itgam'
$ws.Range("D11").Value = 'Inflammatory-Mac'
$ws.Range("E11").Value = [double]"3"
$ws.Range("F11").Value = [double]"1"
$ws.Range("G11").Value = [double]"12.47652966666667"
$ws.Range("H11").Value = [double]"37.429589"
$ws.Range("I11").Value = [double]"0.436533649177391"
$ws.Range("J11").Value = [double]"0.436533649177391"
$ws.Range("K11").Value = [double]"3"
$ws.Range("L11").Value = [double]"1"
$ws.Range("M11").Value = [double]"86.42780700000002"
$ws.Range("N11").Value = [double]"259.283421"
$ws.Range("O11").Value = [double]"0.5839326098770704"
$ws.Range("P11").Value = [double]"0.5839326098770704"
$ws.Range("Q11").Value = [double]"1078.319098060441"
$ws.Range("R11").Value = [double]"9704.87188254397"
$ws.Range("S11").Value = [double]"0.2549062330633153"
$ws.Range("T11").Value = [double]"0.2549062330633153"

# Row 12
$ws.Range("A12").Value = 'Inflammatory-Mac'
$ws.Range("B12").Value = 'Hp'
$ws.Range("C12").Value = 'Itgam'
$ws.Range("D12").Value = 'MuSCs'
$ws.Range("E12").Value = [double]"3"
$ws.Range("F12").Value = [double]"1"
$ws.Range("G12").Value = [double]"12.47652966666667"
$ws.Range("H12").Value = [double]"37.429589"
$ws.Range("I12").Value = [double]"0.436533649177391"
$ws.Range("J12").Value = [double]"0.436533649177391"
$ws.Range("K12").Value = [double]"1"
$ws.Range("L12").Value = [double]"0.3333333333333333"
$ws.Range("M12").Value = [double]"0.006361333333333333"
$ws.Range("N12").Value = [double]"0.019084"
$ws.Range("O12").Value = [double]"4.297910712499435E-05"
$ws.Range("P12").Value = [double]"4.297910712499435E-05"
$ws.Range("Q12").Value = [double]"0.07936736405288888"
$ws.Range("R12").Value = [double]"0.714306276476"
$ws.Range("S12").Value = [double]"1.876182647165979E-05"
$ws.Range("T12").Value = [double]"1.876182647165979E-05"

# Row 13
$ws.Range("A13").Value = 'Inflammatory-Mac'
$ws.Range("B13").Value = 'Hp'
$ws.Range("C13").Value = 'Itgam'
$ws.Range("D13").Value = 'Resolving-Mac'
$ws.Range("E13").Value = [double]"3"
$ws.Range("F13").Value = [double]"1"
$ws.Range("G13").Value = [double]"12.47652966666667"
$ws.Range("H13").Value = [double]"37.429589"
$ws.Range("I13").Value = [double]"0.436533649177391"
$ws.Range("J13").Value = [double]"0.436533649177391"
$ws.Range("K13").Value = [double]"3"
$ws.Range("L13").Value = [double]"1"
$ws.Range("M13").Value = [double]"61.43300833333333"
$ws.Range("N13").Value = [double]"184.299025"
$ws.Range("O13").Value = [double]"0.415060130921558"
$ws.Range("P13").Value = [double]"0.415060130921558"
$ws.Range("Q13").Value = [double]"766.4707509834138"
$ws.Range("R13").Value = [double]"6898.236758850725"
$ws.Range("S13").Value = [double]"0.1811877135792334"
$ws.Range("T13").Value = [double]"0.1811877135792334"

# Row 14
$ws.Range("A14").Value = 'Resolving-Mac'
$ws.Range("B14").Value = 'Hp'
$ws.Range("C14").Value = 'Itgam'
$ws.Range("D14").Value = 'FAPs'
$ws.Range("E14").Value = [double]"1"
$ws.Range("F14").Value = [double]"0.3333333333333333"
$ws.Range("G14").Value = [double]"0.186436"
$ws.Range("H14").Value = [double]"0.559308"
$ws.Range("I14").Value = [double]"0.006523094930433466"
$ws.Range("J14").Value = [double]"0.006523094930433466"
$ws.Range("K14").Value = [double]"1"
$ws.Range("L14").Value = [double]"0.3333333333333333"
$ws.Range("M14").Value = [double]"0.142723"
$ws.Range("N14").Value = [double]"0.428169"
$ws.Range("O14").Value = [double]"0.0009642800942465787"
$ws.Range("P14").Value = [double]"0.0009642800942465787"
$ws.Range("Q14").Value = [double]"0.026608705228"
$ws.Range("R14").Value = [double]"0.239478347052"
$ws.Range("S14").Value = [double]"6.290090594297762E-06"
$ws.Range("T14").Value = [double]"6.290090594297762E-06"

# Row 15
$ws.Range("A15").Value = 'Resolving-Mac'
$ws.Range("B15").Value = 'Hp'
$ws.Range("C15").Value = 'Itgam'
$ws.Range("D15").Value = 'Inflammatory-Mac'
$ws.Range("E15").Value = [double]"1"
$ws.Range("F15").Value = [double]"0.3333333333333333"
$ws.Range("G15").Value = [double]"0.186436"
$ws.Range("H15").Value = [double]"0.559308"
$ws.Range("I15").Value = [double]"0.006523094930433466"
$ws.Range("J15").Value = [double]"0.006523094930433466"
$ws.Range("K15").Value = [double]"3"
$ws.Range("L15").Value = [double]"1"
$ws.Range("M15").Value = [double]"86.42780700000002"
$ws.Range("N15").Value = [double]"259.283421"
$ws.Range("O15").Value = [double]"0.5839326098770704"
$ws.Range("P15").Value = [double]"0.5839326098770704"
$ws.Range("Q15").Value = [double]"16.113254625852"
$ws.Range("R15").Value = [double]"145.019291632668"
$ws.Range("S15").Value = [double]"0.003809047847203901"
$ws.Range("T15").Value = [double]"0.003809047847203901"

# Row 16
$ws.Range("A16").Value = 'Resolving-Mac'
$ws.Range("B16").Value = 'Hp'
$ws.Range("C16").Value = 'Itgam'
$ws.Range("D16").Value = 'MuSCs'
$ws.Range("E16").Value = [double]"1"
$ws.Range("F16").Value = [double]"0.3333333333333333"
$ws.Range("G16").Value = [double]"0.186436"
$ws.Range("H16").Value = [double]"0.559308"
$ws.Range("I16").Value = [double]"0.006523094930433466"
$ws.Range("J16").Value = [double]"0.006523094930433466"
$ws.Range("K16").Value = [double]"1"
$ws.Range("L16").Value = [double]"0.3333333333333333"
$ws.Range("M16").Value = [double]"0.006361333333333333"
$ws.Range("N16").Value = [double]"0.019084"
$ws.Range("O16").Value = [double]"4.297910712499435E-05"
$ws.Range("P16").Value = [double]"4.297910712499435E-05"
$ws.Range("Q16").Value = [double]"0.001185981541333333"
$ws.Range("R16").Value = [double]"0.010673833872"
$ws.Range("S16").Value = [double]"2.803567958016075E-07"
$ws.Range("T16").Value = [double]"2.803567958016075E-07"

# Row 17
$ws.Range("A17").Value = 'Resolving-Mac'
$ws.Range("B17").Value = 'Hp'
$ws.Range("C17").Value = 'Itgam'
$ws.Range("D17").Value = 'Resolving-Mac'
$ws.Range("E17").Value = [double]"1"
$ws.Range("F17").Value = [double]"0.3333333333333333"
$ws.Range("G17").Value = [double]"0.186436"
$ws.Range("H17").Value = [double]"0.559308"
$ws.Range("I17").Value = [double]"0.006523094930433466"
$ws.Range("J17").Value = [double]"0.006523094930433466"
$ws.Range("K17").Value = [double]"3"
$ws.Range("L17").Value = [double]"1"
$ws.Range("M17").Value = [double]"61.43300833333333"
$ws.Range("N17").Value = [double]"184.299025"
$ws.Range("O17").Value = [double]"0.415060130921558"
$ws.Range("P17").Value = [double]"0.415060130921558"
$ws.Range("Q17").Value = [double]"11.45332434163334"
$ws.Range("R17").Value = [double]"103.0799190747"
$ws.Range("S17").Value = [double]"0.002707476635839465"
$ws.Range("T17").Value = [double]"0.002707476635839465"
